# Apply Jove-publication marker layout edit:
#  - Insert HR_Max column; split VR into VR1/VR2
#  - Insert an R2_<name> column after every BRS* column
#  - Refresh row-2 sample values for the re-laid-out columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 1: column headers (final layout, A..BP) ----
$ws.Range("A1").Value = 'Var1'
$ws.Range("B1").Value = 'T_1'
$ws.Range("C1").Value = 'T_2e'
$ws.Range("D1").Value = 'T_2l'
$ws.Range("E1").Value = 'T_3'
$ws.Range("F1").Value = 'T_PRT'
$ws.Range("G1").Value = 'SP_mean'
$ws.Range("H1").Value = 'SP1_Max'
$ws.Range("I1").Value = 'SP2e_min'
$ws.Range("J1").Value = 'SP2l_Max'
$ws.Range("K1").Value = 'SP2l_end'
$ws.Range("L1").Value = 'SP3_end'
$ws.Range("M1").Value = 'SP4e_Max'
$ws.Range("N1").Value = 'HRb_mean'
$ws.Range("O1").Value = 'HRa_mean'
$ws.Range("P1").Value = 'HR2e_min'
$ws.Range("Q1").Value = 'HR2e_end'
$ws.Range("R1").Value = 'HR4e_Max'
$ws.Range("S1").Value = 'HR4e_min'
$ws.Range("T1").Value = 'HR4_min'
$ws.Range("U1").Value = 'HR_Max'
$ws.Range("V1").Value = 'VR1'
$ws.Range("W1").Value = 'VR2'
$ws.Range("X1").Value = 'RRb_mean'
$ws.Range("Y1").Value = 'RRa_mean'
$ws.Range("Z1").Value = 'RR2e_Max'
$ws.Range("AA1").Value = 'RR2e_end'
$ws.Range("AB1").Value = 'RR4e_min'
$ws.Range("AC1").Value = 'RR4e_Max'
$ws.Range("AD1").Value = 'RR4_Max'
$ws.Range("AE1").Value = 'BRSv2eHRTD'
$ws.Range("AF1").Value = 'R2_BRSv2eHRTD'
$ws.Range("AG1").Value = 'BRSv2eRRTD'
$ws.Range("AH1").Value = 'R2_BRSv2eRRTD'
$ws.Range("AI1").Value = 'BRSv2eSPTD'
$ws.Range("AJ1").Value = 'R2_BRSv2eSPTD'
$ws.Range("AK1").Value = 'BRSv2eHRSP'
$ws.Range("AL1").Value = 'R2_BRSv2eHRSP'
$ws.Range("AM1").Value = 'BRSv2eRRSP'
$ws.Range("AN1").Value = 'R2_BRSv2eRRSP'
$ws.Range("AO1").Value = 'BRSa2lSPTD'
$ws.Range("AP1").Value = 'R2_BRSa2lSPTD'
$ws.Range("AQ1").Value = 'BRSv4eHRTD'
$ws.Range("AR1").Value = 'R2_BRSv4eHRTD'
$ws.Range("AS1").Value = 'BRSv4eRRTD'
$ws.Range("AT1").Value = 'R2_BRSv4eRRTD'
$ws.Range("AU1").Value = 'BRSv4eSPTD'
$ws.Range("AV1").Value = 'R2_BRSv4eSPTD'
$ws.Range("AW1").Value = 'BRSv4eHRSP'
$ws.Range("AX1").Value = 'R2_BRSv4eHRSP'
$ws.Range("AY1").Value = 'BRSv4eRRSP'
$ws.Range("AZ1").Value = 'R2_BRSv4eRRSP'
$ws.Range("BA1").Value = 'A'
$ws.Range("BB1").Value = 'B'
$ws.Range("BC1").Value = 'C'
$ws.Range("BD1").Value = 'D'
$ws.Range("BE1").Value = 'E'
$ws.Range("BF1").Value = 'BRSa'
$ws.Range("BG1").Value = 'BRSa1'
$ws.Range("BH1").Value = 'alpha_BRSa'
$ws.Range("BI1").Value = 'beta_BRSa'
$ws.Range("BJ1").Value = 'alpha'
$ws.Range("BK1").Value = 'beta'
$ws.Range("BL1").Value = 'alpha_Area'
$ws.Range("BM1").Value = 'beta_Area'
$ws.Range("BN1").Value = 'BRSa_Area'
$ws.Range("BO1").Value = 'BRSg'
$ws.Range("BP1").Value = 'BRSg1'

# ---- Row 2: sample values (final layout, A..BP) ----
$ws.Range("A2").Value = 'Subject_1'
$ws.Range("B2").Value = 2.6000000000000014
$ws.Range("C2").Value = 5.1999999999999975
$ws.Range("D2").Value = 11.300000000000004
$ws.Range("E2").Value = 1.2999999999999972
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 106.17345476232337
$ws.Range("H2").Value = 119.76829520812441
$ws.Range("I2").Value = 68.433783933674249
$ws.Range("J2").Value = 99.785555523046796
$ws.Range("K2").Value = 99.284195040118391
$ws.Range("L2").Value = 74.178824909778541
$ws.Range("M2").Value = 140.60975609756096
$ws.Range("N2").Value = 87.96411891723038
$ws.Range("O2").Value = 51.356340615340521
$ws.Range("P2").Value = 74.140327865838699
$ws.Range("Q2").Value = 104.54262063419709
$ws.Range("R2").Value = 134.67657342657344
$ws.Range("S2").Value = 68.592657342657361
$ws.Range("T2").Value = 40.788579197824703
$ws.Range("U2").Value = 134.52435061989851
$ws.Range("V2").Value = 3.3018206585081513
$ws.Range("W2").Value = 3.2980886626978378
$ws.Range("X2").Value = 0.68209629947475348
$ws.Range("Y2").Value = 1.1683075406287331
$ws.Range("Z2").Value = 809.27616220653272
$ws.Range("AA2").Value = 573.92860094778712
$ws.Range("AB2").Value = 445.51178035957679
$ws.Range("AC2").Value = 874.72919587103331
$ws.Range("AD2").Value = 1470.9999999999966
$ws.Range("AE2").Value = 6.2491899886815476
$ws.Range("AF2").Value = 0.99670035961287151
$ws.Range("AG2").Value = -48.490762596464826
$ws.Range("AH2").Value = 0.98899643535971915
$ws.Range("AI2").Value = -10.183422893894138
$ws.Range("AJ2").Value = 0.91126954647069136
$ws.Range("AK2").Value = -0.55910811347111355
$ws.Range("AL2").Value = 0.9079239014293693
$ws.Range("AM2").Value = 4.4591120532952324
$ws.Range("AN2").Value = 0.95173090383434988
$ws.Range("AO2").Value = 3.3116017471065144
$ws.Range("AP2").Value = 0.98485773910417751
$ws.Range("AQ2").Value = -87.643001924327436
$ws.Range("AR2").Value = 0.92272272397380195
$ws.Range("AS2").Value = 511.12824459399485
$ws.Range("AT2").Value = 0.85508285118565008
$ws.Range("AU2").Value = 25.67321372693884
$ws.Range("AV2").Value = 0.98142924898858108
$ws.Range("AW2").Value = -4.0557656914351696
$ws.Range("AX2").Value = 0.93192449338569505
$ws.Range("AY2").Value = 23.762742432444526
$ws.Range("AZ2").Value = 0.87164539884357051
$ws.Range("BA2").Value = -37.739670828649125
$ws.Range("BB2").Value = -25.10537013033985
$ws.Range("BC2").Value = 30.850411106444142
$ws.Range("BD2").Value = 31.994629852544833
$ws.Range("BE2").Value = 34.436301335237587
$ws.Range("BF2").Value = 25.10537013033985
$ws.Range("BG2").Value = 56.568698426404012
$ws.Range("BH2").Value = 2.7301248766764719
$ws.Range("BI2").Value = 31.994629852544833
$ws.Range("BJ2").Value = 73.197328513418512
$ws.Range("BK2").Value = 88.363549648812253
$ws.Range("BL2").Value = 196.7257017868624
$ws.Range("BM2").Value = 17.501355549123332
$ws.Range("BN2").Value = 3442.9664526226861
$ws.Range("BO2").Value = 0.11194765855063651
$ws.Range("BP2").Value = 0.25224616499240116

# ---- Column widths to match final layout (best effort; Excel quantizes to pixels) ----
$ws.Columns("A").ColumnWidth = 4.307291666666667
$ws.Columns("B").ColumnWidth = 3.3072916666666665
$ws.Columns("C").ColumnWidth = 4.451822916666667
$ws.Columns("D").ColumnWidth = 3.8776041666666665
$ws.Columns("E").ColumnWidth = 3.3072916666666665
$ws.Columns("F").ColumnWidth = 5.592447916666667
$ws.Columns("G").ColumnWidth = 8.451822916666666
$ws.Columns("H").ColumnWidth = 8.166666666666666
$ws.Columns("I").ColumnWidth = 9.022135416666666
$ws.Columns("J").ColumnWidth = 8.736979166666666
$ws.Columns("K").ColumnWidth = 8.451822916666666
$ws.Columns("L").ColumnWidth = 7.877604166666667
$ws.Columns("M").ColumnWidth = 9.307291666666666
$ws.Columns("N").ColumnWidth = 9.877604166666666
$ws.Columns("O").ColumnWidth = 9.736979166666666
$ws.Columns("P").ColumnWidth = 9.307291666666666
$ws.Columns("Q").ColumnWidth = 9.307291666666666
$ws.Columns("R").ColumnWidth = 9.592447916666666
$ws.Columns("S").ColumnWidth = 9.307291666666666
$ws.Columns("T").ColumnWidth = 8.166666666666666
$ws.Columns("U").ColumnWidth = 7.451822916666667
$ws.Columns("V").ColumnWidth = 3.7369791666666665
$ws.Columns("W").ColumnWidth = 3.7369791666666665
$ws.Columns("X").ColumnWidth = 9.736979166666666
$ws.Columns("Y").ColumnWidth = 9.592447916666666
$ws.Columns("Z").ColumnWidth = 9.451822916666666
$ws.Columns("AA").ColumnWidth = 9.166666666666666
$ws.Columns("AB").ColumnWidth = 9.166666666666666
$ws.Columns("AC").ColumnWidth = 9.451822916666666
$ws.Columns("AD").ColumnWidth = 8.307291666666666
$ws.Columns("AE").ColumnWidth = 11.451822916666666
$ws.Columns("AF").ColumnWidth = 14.592447916666666
$ws.Columns("AG").ColumnWidth = 11.307291666666666
$ws.Columns("AH").ColumnWidth = 14.451822916666666
$ws.Columns("AI").ColumnWidth = 11.166666666666666
$ws.Columns("AJ").ColumnWidth = 14.307291666666666
$ws.Columns("AK").ColumnWidth = 11.307291666666666
$ws.Columns("AL").ColumnWidth = 14.451822916666666
$ws.Columns("AM").ColumnWidth = 11.166666666666666
$ws.Columns("AN").ColumnWidth = 14.307291666666666
$ws.Columns("AO").ColumnWidth = 10.592447916666666
$ws.Columns("AP").ColumnWidth = 13.736979166666666
$ws.Columns("AQ").ColumnWidth = 11.451822916666666
$ws.Columns("AR").ColumnWidth = 14.592447916666666
$ws.Columns("AS").ColumnWidth = 11.307291666666666
$ws.Columns("AT").ColumnWidth = 14.451822916666666
$ws.Columns("AU").ColumnWidth = 11.166666666666666
$ws.Columns("AV").ColumnWidth = 14.307291666666666
$ws.Columns("AW").ColumnWidth = 11.307291666666666
$ws.Columns("AX").ColumnWidth = 14.451822916666666
$ws.Columns("AY").ColumnWidth = 11.166666666666666
$ws.Columns("AZ").ColumnWidth = 14.307291666666666
$ws.Columns("BA").ColumnWidth = 1.5924479166666665
$ws.Columns("BB").ColumnWidth = 1.4518229166666665
$ws.Columns("BC").ColumnWidth = 1.4518229166666665
$ws.Columns("BD").ColumnWidth = 1.5924479166666665
$ws.Columns("BE").ColumnWidth = 1.3072916666666665
$ws.Columns("BF").ColumnWidth = 4.592447916666667
$ws.Columns("BG").ColumnWidth = 5.592447916666667
$ws.Columns("BH").ColumnWidth = 10.451822916666666
$ws.Columns("BI").ColumnWidth = 9.592447916666666
$ws.Columns("BJ").ColumnWidth = 5.166666666666667
$ws.Columns("BK").ColumnWidth = 4.307291666666667
$ws.Columns("BL").ColumnWidth = 10.307291666666666
$ws.Columns("BM").ColumnWidth = 9.451822916666666
$ws.Columns("BN").ColumnWidth = 9.736979166666666
$ws.Columns("BO").ColumnWidth = 4.592447916666667
$ws.Columns("BP").ColumnWidth = 5.592447916666667
